# Slide 6 ("QUESTION 4") -> "TextBox 8" (shape id 9) holding the
# "Metric Used to Evaluate..." explanation. We append two new paragraphs
# and grow the textbox to its new size.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(13)

$tr = $sh.TextFrame.TextRange
[void]$tr.InsertAfter("`rHence, CA has the best Conversion Rate of Clicks, `rmaking it the most efficient store ")

# Resize the textbox to match the new (auto-fit) extent recorded in the
# target deck: cx="4842864" cy="2585323" EMU -> points (914400 EMU/in, 72 pt/in)
$sh.Width = 4842864 / 914400 * 72
$sh.Height = 2585323 / 914400 * 72
